# #56 Created Method and View to list frinds of user
#
# - Adds a new Controller "method" row (row 13, left table) for the
#   new FrindsList method that returns the list of a user's friends.
# - Updates the existing "FriendsList" View row (row 14, right table)
#   text to reflect that the page/method are now implemented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: "FriendsList" view - refreshed notes now that it's done ---
$ws.Range("K14").Value = "Pagina que mostra todos os amigos do user"
$ws.Range("L14").Value = "Recebe uma Lista de Users"
$ws.Range("M14").Value = "X"
$ws.Range("N14").Value = "X"
$ws.Rows.Item(14).AutoFit()

# --- Row 13: new Controller method "FrindsList" ------------------------
$ws.Range("B13").Value = "User"
$ws.Range("C13").Value = "FrindsList"
$ws.Range("D13").Value = "Pesquisa todos os users amigos do user `nid recebido por parametro"
$ws.Range("D13").WrapText = $true
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(13).RowHeight = 30

# --- View/selection state, matches the author's saved window position --
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("N10").Select()
